$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename: "Athlete Class" -> "Class" ---
$ws.Range("A1").Value = "Class"

# --- New "Match" state-machine rows (32-40) ---
$ws.Range("A32").Value = "Match"
$ws.Range("C32").Value = "MatchStartState"
$ws.Range("D32").Value = "Awake"

$ws.Range("A33").Value = "Match"
$ws.Range("B33").Value = "MatchStartState"
$ws.Range("C33").Value = "PrePointState"
$ws.Range("D33").Value = "MatchInfo.Initialize() [event]"

$ws.Range("A34").Value = "Match"
$ws.Range("B34").Value = "PrePointState"
$ws.Range("C34").Value = "ServeState"
$ws.Range("D34").Value = "OnInteract"

$ws.Range("A35").Value = "Match"
$ws.Range("B35").Value = "ServeState"
$ws.Range("C35").Value = "InPlayState"
$ws.Range("D35").Value = "OnBallServed"

$ws.Range("A36").Value = "Match"
$ws.Range("B36").Value = "InPlayState"
$ws.Range("C36").Value = "postPointState"
$ws.Range("D36").Value = "OnBallHitGround"

$ws.Range("A37").Value = "Match"
$ws.Range("B37").Value = "postPointState"
$ws.Range("C37").Value = "MatchEndState"
$ws.Range("D37").Value = "GameEnd == true"

$ws.Range("A38").Value = "Match"
$ws.Range("B38").Value = "postPointState"
$ws.Range("C38").Value = "PrePointState"
$ws.Range("D38").Value = "timeUntilChange [event]"

$ws.Range("A39").Value = "Match"
$ws.Range("B39").Value = "any"
$ws.Range("C39").Value = "PauseState"
$ws.Range("D39").Value = "OnPause"

$ws.Range("A40").Value = "Match"
$ws.Range("B40").Value = "PauseState"
$ws.Range("C40").Value = "variable"
$ws.Range("D40").Value = "OnUnPause"

# --- Italicize the "any" wildcard cells (B8, B39, C40) ---
# B39 starts from the plain (color-bearing) default font, so italicizing it here
# mints the new italic font/style cleanly; then we fan that exact style out to
# the other two "any" cells via a formats-only paste so no duplicate styles are
# minted.
$ws.Range("B39").Font.Italic = $true
$ws.Range("B39").Copy()
$ws.Range("C40").PasteSpecial(-4122)
$ws.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Final selection state ---
[void]$ws.Range("D41").Select()
